$wb = $excel.ActiveWorkbook

# ---- Sheet 'Recommandations' ----
$ws1 = $wb.Worksheets.Item("Recommandations")

# Remove the last data row (old row 48) so the table shrinks from 47 data rows to 46
$ws1.Range("A48:G48").EntireRow.Delete()

$data1 = New-Object 'object[,]' 46,7
$data1[0,0] = 'BRVM - SERVICES PUBLICS'
$data1[0,1] = 0
$data1[0,2] = 8
$data1[0,3] = 3356.33
$data1[0,4] = 105.88
$data1[0,5] = '🟡 Observer'
$data1[0,6] = '➖ Neutre'
$data1[1,0] = 'AIR LIQUIDE CI'
$data1[1,1] = 0
$data1[1,2] = 4
$data1[1,3] = 2805
$data1[1,4] = 715
$data1[1,5] = '🟡 Observer'
$data1[1,6] = '➖ Neutre'
$data1[2,0] = 'NEI-CEDA CI'
$data1[2,1] = 0
$data1[2,2] = 4
$data1[2,3] = 2775
$data1[2,4] = 725
$data1[2,5] = '🟡 Observer'
$data1[2,6] = '➖ Neutre'
$data1[3,0] = 'BRVM - AUTRES SECTEURS'
$data1[3,1] = 0
$data1[3,2] = 4
$data1[3,3] = 2481.33
$data1[3,4] = 614.21
$data1[3,5] = '🟡 Observer'
$data1[3,6] = '➖ Neutre'
$data1[4,0] = 'BRVM - DISTRIBUTION'
$data1[4,1] = 0
$data1[4,2] = 4
$data1[4,3] = 2059.92
$data1[4,4] = 489.13
$data1[4,5] = '🟡 Observer'
$data1[4,6] = '➖ Neutre'
$data1[5,0] = 'BRVM - AGRICULTURE'
$data1[5,1] = 0
$data1[5,2] = 4
$data1[5,3] = 1499.13
$data1[5,4] = 372.93
$data1[5,5] = '🟡 Observer'
$data1[5,6] = '➖ Neutre'
$data1[6,0] = 'BRVM - TRANSPORT'
$data1[6,1] = 0
$data1[6,2] = 4
$data1[6,3] = 1428.25
$data1[6,4] = 357.37
$data1[6,5] = '🟡 Observer'
$data1[6,6] = '➖ Neutre'
$data1[7,0] = 'BRVM - CONSOMMATION DISCRETIONNAIRE'
$data1[7,1] = 0
$data1[7,2] = 4
$data1[7,3] = 731.27
$data1[7,4] = 172.56
$data1[7,5] = '🟡 Observer'
$data1[7,6] = '➖ Neutre'
$data1[8,0] = 'BRVM - FINANCES'
$data1[8,1] = 0
$data1[8,2] = 4
$data1[8,3] = 560.5599999999999
$data1[8,4] = 140.3
$data1[8,5] = '🟡 Observer'
$data1[8,6] = '➖ Neutre'
$data1[9,0] = 'BRVM-PRESTIGE'
$data1[9,1] = 0
$data1[9,2] = 4
$data1[9,3] = 559.34
$data1[9,4] = 139.34
$data1[9,5] = '🟡 Observer'
$data1[9,6] = '➖ Neutre'
$data1[10,0] = 'BRVM - SERVICES FINANCIERS'
$data1[10,1] = 0
$data1[10,2] = 4
$data1[10,3] = 550.91
$data1[10,4] = 137.89
$data1[10,5] = '🟡 Observer'
$data1[10,6] = '➖ Neutre'
$data1[11,0] = 'BRVM - INDUSTRIELS'
$data1[11,1] = 0
$data1[11,2] = 4
$data1[11,3] = 519.2
$data1[11,4] = 128.07
$data1[11,5] = '🟡 Observer'
$data1[11,6] = '➖ Neutre'
$data1[12,0] = 'BRVM - ENERGIE'
$data1[12,1] = 0
$data1[12,2] = 4
$data1[12,3] = 441.47
$data1[12,4] = 108.05
$data1[12,5] = '🟡 Observer'
$data1[12,6] = '➖ Neutre'
$data1[13,0] = 'BRVM-PRINCIPAL                    (**)'
$data1[13,1] = 0
$data1[13,2] = 2
$data1[13,3] = 410.25
$data1[13,4] = 204.33
$data1[13,5] = '🟡 Observer'
$data1[13,6] = '➖ Neutre'
$data1[14,0] = 'BRVM - TELECOMMUNICATIONS'
$data1[14,1] = 0
$data1[14,2] = 4
$data1[14,3] = 381.27
$data1[14,4] = 95.13
$data1[14,5] = '🟡 Observer'
$data1[14,6] = '➖ Neutre'
$data1[15,0] = 'BRVM - INDUSTRIE                 (**)'
$data1[15,1] = 0
$data1[15,2] = 1
$data1[15,3] = 219.23
$data1[15,4] = 219.23
$data1[15,5] = '🟡 Observer'
$data1[15,6] = '➖ Neutre'
$data1[16,0] = 'BRVM - INDUSTRIE                  (**)'
$data1[16,1] = 0
$data1[16,2] = 1
$data1[16,3] = 218.66
$data1[16,4] = 218.66
$data1[16,5] = '🟡 Observer'
$data1[16,6] = '➖ Neutre'
$data1[17,0] = 'BRVM - INDUSTRIE              (**)'
$data1[17,1] = 0
$data1[17,2] = 1
$data1[17,3] = 214.53
$data1[17,4] = 214.53
$data1[17,5] = '🟡 Observer'
$data1[17,6] = '➖ Neutre'
$data1[18,0] = 'BRVM-PRINCIPAL                (**)'
$data1[18,1] = 0
$data1[18,2] = 1
$data1[18,3] = 205.03
$data1[18,4] = 205.03
$data1[18,5] = '🟡 Observer'
$data1[18,6] = '➖ Neutre'
$data1[19,0] = 'BRVM - CONSOMMATION DE BASE        (**)'
$data1[19,1] = 0
$data1[19,2] = 1
$data1[19,3] = 194.82
$data1[19,4] = 194.82
$data1[19,5] = '🟡 Observer'
$data1[19,6] = '➖ Neutre'
$data1[20,0] = 'BRVM - CONSOMMATION DE BASE          (**)'
$data1[20,1] = 0
$data1[20,2] = 1
$data1[20,3] = 193.9
$data1[20,4] = 193.9
$data1[20,5] = '🟡 Observer'
$data1[20,6] = '➖ Neutre'
$data1[21,0] = 'BRVM - CONSOMMATION DE BASE               (**)'
$data1[21,1] = 0
$data1[21,2] = 1
$data1[21,3] = 191.47
$data1[21,4] = 191.47
$data1[21,5] = '🟡 Observer'
$data1[21,6] = '➖ Neutre'
$data1[22,0] = 'SERVAIR ABIDJAN CI (ABJC)'
$data1[22,1] = 4
$data1[22,2] = 0
$data1[22,3] = 25.86
$data1[22,4] = 3.75
$data1[22,5] = '🟢 Achat'
$data1[22,6] = '✅ Renforcer'
$data1[23,0] = 'UNIWAX CI (UNXC)'
$data1[23,1] = 2
$data1[23,2] = 0
$data1[23,3] = 14.63
$data1[23,4] = 7.21
$data1[23,5] = '🟡 Observer'
$data1[23,6] = '➖ Neutre'
$data1[24,0] = 'SOLIBRA CI (SLBC)'
$data1[24,1] = 1
$data1[24,2] = 0
$data1[24,3] = 7.48
$data1[24,4] = 7.48
$data1[24,5] = '🟡 Observer'
$data1[24,6] = '➖ Neutre'
$data1[25,0] = 'VIVO ENERGY CI (SHEC)'
$data1[25,1] = 1
$data1[25,2] = 0
$data1[25,3] = 7
$data1[25,4] = 7
$data1[25,5] = '🟡 Observer'
$data1[25,6] = '➖ Neutre'
$data1[26,0] = 'TRACTAFRIC MOTORS CI (PRSC)'
$data1[26,1] = 2
$data1[26,2] = 0
$data1[26,3] = 5.35
$data1[26,4] = 2.29
$data1[26,5] = '🟡 Observer'
$data1[26,6] = '➖ Neutre'
$data1[27,0] = 'ECOBANK TRANS. INCORP. TG (ETIT)'
$data1[27,1] = 1
$data1[27,2] = 0
$data1[27,3] = 5
$data1[27,4] = 5
$data1[27,5] = '🟡 Observer'
$data1[27,6] = '➖ Neutre'
$data1[28,0] = 'SOGB CI (SOGC)'
$data1[28,1] = 1
$data1[28,2] = 0
$data1[28,3] = 2.07
$data1[28,4] = 2.07
$data1[28,5] = '🟡 Observer'
$data1[28,6] = '➖ Neutre'
$data1[29,0] = 'NEI-CEDA CI (NEIC)'
$data1[29,1] = 1
$data1[29,2] = 0
$data1[29,3] = 1.45
$data1[29,4] = 1.45
$data1[29,5] = '🟡 Observer'
$data1[29,6] = '➖ Neutre'
$data1[30,0] = 'SAPH CI (SPHC)'
$data1[30,1] = 1
$data1[30,2] = 0
$data1[30,3] = 1.44
$data1[30,4] = 1.44
$data1[30,5] = '🟡 Observer'
$data1[30,6] = '➖ Neutre'
$data1[31,0] = 'SICOR CI (SICC)'
$data1[31,1] = 1
$data1[31,2] = 1
$data1[31,3] = 0.55
$data1[31,4] = 7.46
$data1[31,5] = '🟡 Observer'
$data1[31,6] = '👀 À surveiller'
$data1[32,0] = 'SAFCA CI (SAFC)'
$data1[32,1] = 2
$data1[32,2] = 2
$data1[32,3] = 0.01
$data1[32,4] = -7.49
$data1[32,5] = '🟡 Observer'
$data1[32,6] = '👀 À surveiller'
$data1[33,0] = 'TOTAL'
$data1[33,1] = 0
$data1[33,2] = 4
$data1[33,3] = 0
$data1[33,4] = 0
$data1[33,5] = '🟡 Observer'
$data1[33,6] = '➖ Neutre'
$data1[34,0] = 'ONATEL BF (ONTBF)'
$data1[34,1] = 1
$data1[34,2] = 1
$data1[34,3] = -0.65
$data1[34,4] = 5.49
$data1[34,5] = '🟡 Observer'
$data1[34,6] = '👀 À surveiller'
$data1[35,0] = 'BANK OF AFRICA ML (BOAM)'
$data1[35,1] = 0
$data1[35,2] = 1
$data1[35,3] = -2.19
$data1[35,4] = -2.19
$data1[35,5] = '🟡 Observer'
$data1[35,6] = '➖ Neutre'
$data1[36,0] = 'SICABLE CI (CABC)'
$data1[36,1] = 1
$data1[36,2] = 2
$data1[36,3] = -2.66
$data1[36,4] = -7.37
$data1[36,5] = '🟡 Observer'
$data1[36,6] = '➖ Neutre'
$data1[37,0] = 'FILTISAC CI (FTSC)'
$data1[37,1] = 0
$data1[37,2] = 1
$data1[37,3] = -3.07
$data1[37,4] = -3.07
$data1[37,5] = '🟡 Observer'
$data1[37,6] = '➖ Neutre'
$data1[38,0] = 'AIR LIQUIDE CI (SIVC)'
$data1[38,1] = 0
$data1[38,2] = 1
$data1[38,3] = -3.5
$data1[38,4] = -3.5
$data1[38,5] = '🟡 Observer'
$data1[38,6] = '➖ Neutre'
$data1[39,0] = 'ECOBANK COTE D''''IVOIRE (ECOC)'
$data1[39,1] = 0
$data1[39,2] = 1
$data1[39,3] = -3.85
$data1[39,4] = -3.85
$data1[39,5] = '🟡 Observer'
$data1[39,6] = '➖ Neutre'
$data1[40,0] = 'SOCIETE GENERALE COTE D''IVOIRE (SGBC)'
$data1[40,1] = 0
$data1[40,2] = 1
$data1[40,3] = -3.91
$data1[40,4] = -3.91
$data1[40,5] = '🟡 Observer'
$data1[40,6] = '➖ Neutre'
$data1[41,0] = 'BERNABE CI (BNBC)'
$data1[41,1] = 0
$data1[41,2] = 1
$data1[41,3] = -5.66
$data1[41,4] = -5.66
$data1[41,5] = '🟡 Observer'
$data1[41,6] = '➖ Neutre'
$data1[42,0] = 'CIE CI (CIEC)'
$data1[42,1] = 0
$data1[42,2] = 1
$data1[42,3] = -6.3
$data1[42,4] = -6.3
$data1[42,5] = '🟡 Observer'
$data1[42,6] = '➖ Neutre'
$data1[43,0] = 'SETAO CI (STAC)'
$data1[43,1] = 0
$data1[43,2] = 2
$data1[43,3] = -10.33
$data1[43,4] = -7.08
$data1[43,5] = '🟡 Observer'
$data1[43,6] = '➖ Neutre'
$data1[44,0] = 'CFAO MOTORS CI (CFAC)'
$data1[44,1] = 1
$data1[44,2] = 3
$data1[44,3] = -13.2
$data1[44,4] = -5.75
$data1[44,5] = '🔴 Vente'
$data1[44,6] = '⚠️ Risque de décrochage'
$data1[45,0] = 'TOTALENERGIES MARKETING CI (TTLC)'
$data1[45,1] = 0
$data1[45,2] = 2
$data1[45,3] = -14.35
$data1[45,4] = -7.35
$data1[45,5] = '🟡 Observer'
$data1[45,6] = '➖ Neutre'
$ws1.Range("A2:G47").Value = $data1

# ---- Sheet 'Top_YTD' ----
$ws2 = $wb.Worksheets.Item("Top_YTD")
$data2 = New-Object 'object[,]' 10,2
$data2[0,0] = 'BRVM - SERVICES PUBLICS'
$data2[0,1] = 8901071.23
$data2[1,0] = 'AIR LIQUIDE CI'
$data2[1,1] = 411572.8
$data2[2,0] = 'NEI-CEDA CI'
$data2[2,1] = 396296
$data2[3,0] = 'BRVM - AUTRES SECTEURS'
$data2[3,1] = 269117.58
$data2[4,0] = 'BRVM - DISTRIBUTION'
$data2[4,1] = 142718.36
$data2[5,0] = 'BRVM - AGRICULTURE'
$data2[5,1] = 50711.9
$data2[6,0] = 'BRVM - TRANSPORT'
$data2[6,1] = 43537.29
$data2[7,0] = 'BRVM - CONSOMMATION DISCRETIONNAIRE'
$data2[7,1] = 6287.63
$data2[8,0] = 'BRVM - FINANCES'
$data2[8,1] = 3225.42
$data2[9,0] = 'BRVM-PRESTIGE'
$data2[9,1] = 3208.46
$ws2.Range("A2:B11").Value = $data2

